$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B-G (Component Analysis filtering update)
$data = @(
    @(2,  0.001974836808913309, 0.7564986997685976, 1.331115891817462, 1.153739958490414, 1.165451647513039, 50),
    @(3,  0.1715060524283095,   0.7455020267087988, 1.556349277421864, 1.247537284982643, 1.248768565063892, 48),
    @(4,  -0.1210599245473831,  0.6932517703744887, 1.42325916890264,  1.19300426189626,  1.199405667105797, 48),
    @(5,  0.1732007841850449,   0.7115279881092994, 1.411811716246263, 1.188196833965763, 1.188494915725877, 46),
    @(6,  0.03614660433360007,  0.5611928274455102, 0.7065544791159908, 0.8405679503264389, 0.849279849945092, 45),
    @(7,  -0.03741410742334574, 0.6726782106122756, 1.306474751485271, 1.143011264811188, 1.159578656834397, 34),
    @(8,  -0.04498115623876454, 0.6806687468339884, 1.217036297618114, 1.103193680918321, 1.119366851947603, 33),
    @(9,  0.008498424339290329, 0.6040240707239352, 0.6504015573506643, 0.8064747716765009, 0.832877315741465, 16),
    @(10, -0.005607768909606597, 0.504936332959623,  0.4301487342848143, 0.6558572514540144, 0.6913089735443678, 10),
    @(11, -0.174551384566527,   0.5019007433328729, 0.3135597735887778, 0.5599640824095575, 0.5948650978737597, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
